$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.76993966666667
$ws.Range("H2").Value = 56.309819
$ws.Range("I2").Value = 0.1007685501185251
$ws.Range("J2").Value = 0.1007685501185251
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 71.23553112006256
$ws.Range("R2").Value = 641.119780080563
$ws.Range("S2").Value = 0.001052239061769869
$ws.Range("T2").Value = 0.001052239061769869
$ws.Range("G3").Value = 18.76993966666667
$ws.Range("H3").Value = 56.309819
$ws.Range("I3").Value = 0.1007685501185251
$ws.Range("J3").Value = 0.1007685501185251
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 4568.160319263948
$ws.Range("R3").Value = 41113.44287337553
$ws.Range("S3").Value = 0.06747751652549788
$ws.Range("T3").Value = 0.06747751652549788
$ws.Range("G4").Value = 18.76993966666667
$ws.Range("H4").Value = 56.309819
$ws.Range("I4").Value = 0.1007685501185251
$ws.Range("J4").Value = 0.1007685501185251
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 559.3699356539496
$ws.Range("R4").Value = 5034.329420885547
$ws.Range("S4").Value = 0.008262602763258051
$ws.Range("T4").Value = 0.00826260276325805
$ws.Range("G5").Value = 18.76993966666667
$ws.Range("H5").Value = 56.309819
$ws.Range("I5").Value = 0.1007685501185251
$ws.Range("J5").Value = 0.1007685501185251
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 1623.164180920176
$ws.Range("R5").Value = 14608.47762828158
$ws.Range("S5").Value = 0.02397619176799932
$ws.Range("T5").Value = 0.02397619176799932
$ws.Range("I6").Value = 0.5130361557055731
$ws.Range("J6").Value = 0.5130361557055731
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 362.6766783137728
$ws.Range("R6").Value = 3264.090104823954
$ws.Range("S6").Value = 0.005357194109657135
$ws.Range("T6").Value = 0.005357194109657135
$ws.Range("I7").Value = 0.5130361557055731
$ws.Range("J7").Value = 0.5130361557055731
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("S7").Value = 0.3435437508437121
$ws.Range("T7").Value = 0.3435437508437121
$ws.Range("I8").Value = 0.5130361557055731
$ws.Range("J8").Value = 0.5130361557055731
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 2847.882608885714
$ws.Range("R8").Value = 25630.94347997142
$ws.Range("S8").Value = 0.04206683486859918
$ws.Range("T8").Value = 0.04206683486859917
$ws.Range("I9").Value = 0.5130361557055731
$ws.Range("J9").Value = 0.5130361557055731
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 8263.906848702218
$ws.Range("R9").Value = 74375.16163831996
$ws.Range("S9").Value = 0.1220683758836047
$ws.Range("T9").Value = 0.1220683758836047
$ws.Range("G10").Value = 20.061603
$ws.Range("H10").Value = 60.184809
$ws.Range("I10").Value = 0.1077029912330274
$ws.Range("J10").Value = 0.1077029912330274
$ws.Range("M10").Value = 3.795192333333334
$ws.Range("N10").Value = 11.385577
$ws.Range("O10").Value = 0.01044213755712683
$ws.Range("P10").Value = 0.01044213755712683
$ws.Range("Q10").Value = 76.13764189997701
$ws.Range("R10").Value = 685.2387770997931
$ws.Range("S10").Value = 0.001124649449769298
$ws.Range("T10").Value = 0.001124649449769298
$ws.Range("G11").Value = 20.061603
$ws.Range("H11").Value = 60.184809
$ws.Range("I11").Value = 0.1077029912330274
$ws.Range("J11").Value = 0.1077029912330274
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.6696287328350964
$ws.Range("P11").Value = 0.6696287328350964
$ws.Range("Q11").Value = 4882.520689620397
$ws.Range("R11").Value = 43942.68620658357
$ws.Range("S11").Value = 0.07212101754192166
$ws.Range("T11").Value = 0.07212101754192166
$ws.Range("G12").Value = 20.061603
$ws.Range("H12").Value = 60.184809
$ws.Range("I12").Value = 0.1077029912330274
$ws.Range("J12").Value = 0.1077029912330274
$ws.Range("M12").Value = 29.801371
$ws.Range("N12").Value = 89.404113
$ws.Range("O12").Value = 0.08199584844219236
$ws.Range("P12").Value = 0.08199584844219235
$ws.Range("Q12").Value = 597.863273857713
$ws.Range("R12").Value = 5380.769464719417
$ws.Range("S12").Value = 0.00883119814591409
$ws.Range("T12").Value = 0.00883119814591409
$ws.Range("G13").Value = 20.061603
$ws.Range("H13").Value = 60.184809
$ws.Range("I13").Value = 0.1077029912330274
$ws.Range("J13").Value = 0.1077029912330274
$ws.Range("M13").Value = 86.47679266666667
$ws.Range("N13").Value = 259.430378
$ws.Range("O13").Value = 0.2379332811655844
$ws.Range("P13").Value = 0.2379332811655844
$ws.Range("Q13").Value = 1734.863083191978
$ws.Range("R13").Value = 15613.7677487278
$ws.Range("S13").Value = 0.02562612609542239
$ws.Range("T13").Value = 0.02562612609542239
$ws.Range("G14").Value = 51.87415833333333
$ws.Range("H14").Value = 155.622475
$ws.Range("I14").Value = 0.2784923029428744
$ws.Range("J14").Value = 0.2784923029428744
$ws.Range("M14").Value = 3.795192333333334
$ws.Range("N14").Value = 11.385577
$ws.Range("O14").Value = 0.01044213755712683
$ws.Range("P14").Value = 0.01044213755712683
$ws.Range("Q14").Value = 196.8724080047861
$ws.Range("R14").Value = 1771.851672043075
$ws.Range("S14").Value = 0.002908054935930532
$ws.Range("T14").Value = 0.002908054935930532
$ws.Range("G15").Value = 51.87415833333333
$ws.Range("H15").Value = 155.622475
$ws.Range("I15").Value = 0.2784923029428744
$ws.Range("J15").Value = 0.2784923029428744
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.6696287328350964
$ws.Range("P15").Value = 0.6696287328350964
$ws.Range("Q15").Value = 12624.94583903112
$ws.Range("R15").Value = 113624.5125512801
$ws.Range("S15").Value = 0.1864864479239647
$ws.Range("T15").Value = 0.1864864479239647
$ws.Range("G16").Value = 51.87415833333333
$ws.Range("H16").Value = 155.622475
$ws.Range("I16").Value = 0.2784923029428744
$ws.Range("J16").Value = 0.2784923029428744
$ws.Range("M16").Value = 29.801371
$ws.Range("N16").Value = 89.404113
$ws.Range("O16").Value = 0.08199584844219236
$ws.Range("P16").Value = 0.08199584844219235
$ws.Range("Q16").Value = 1545.921037804408
$ws.Range("R16").Value = 13913.28934023967
$ws.Range("S16").Value = 0.02283521266442105
$ws.Range("T16").Value = 0.02283521266442104
$ws.Range("G17").Value = 51.87415833333333
$ws.Range("H17").Value = 155.622475
$ws.Range("I17").Value = 0.2784923029428744
$ws.Range("J17").Value = 0.2784923029428744
$ws.Range("M17").Value = 86.47679266666667
$ws.Range("N17").Value = 259.430378
$ws.Range("O17").Value = 0.2379332811655844
$ws.Range("P17").Value = 0.2379332811655844
$ws.Range("Q17").Value = 4485.910834949505
$ws.Range("R17").Value = 40373.19751454556
$ws.Range("S17").Value = 0.06626258741855802
$ws.Range("T17").Value = 0.06626258741855802
